# Applies the tracked-changes cleanup described in the commit:
#   - Relocates the "_GoBack" bookmark from the (now empty) paragraph
#     just before the end of the document to the very start of the
#     document's first paragraph (right after its <w:pPr>, before the
#     first run) -- this is what Word does when the last edit position
#     was at the top of the document.
#   - Collapses several runs that were split only because of how the
#     previous edit session touched them back into single runs holding
#     the full, unbroken sentence (pure run-coalescing, no text-content
#     change).
#
# All of this is done purely through Find/Replace + Bookmarks, i.e. the
# same primitives a human editor driving Word would use; no raw XML
# surgery.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark to the document start.
#
# Word marks the position of the most recent edit with the hidden
# "_GoBack" bookmark. Re-adding it at Range(0,0) reproduces that. A
# bookmark anchored on a truly empty (0,0) range at the very start of
# the story sometimes gets its end marker pushed into the *next*
# paragraph by naive range math, so we anchor it on a transient
# one-character range instead (insert a placeholder, bookmark that
# character, then delete it again) -- this keeps bookmarkStart/bookmarkEnd
# adjacent, right where Word puts them.
# ---------------------------------------------------------------------
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$placeholderRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$d.Range(0, 1).Text = ""

# ---------------------------------------------------------------------
# 2) Re-join the runs that make up each sentence. A Find/Replace whose
# match text is identical to itself (same old/new text) is a no-op on
# content but makes Word re-segment the run(s) it touches, merging them
# with immediately-adjacent runs that already share identical
# formatting. That is exactly the "de-fragmentation" shown in the diff.
# ---------------------------------------------------------------------

function Coalesce-Text([string]$text) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 1) | Out-Null
}

# "Было о" + "чень трудно ..." -> "Было очень трудно ..."
Coalesce-Text "Было о"

# "...в них пуль," + " я использовал сначала " -> one run
Coalesce-Text "пуль,"

# "метод " + "." -> "метод ."
Coalesce-Text "метод "

# "()" + ", но в данном случае ... с помощью " (+ the quoted-word runs in
# between) -> one run; this same edit also re-joins the later
# "метода " + "." pair further down the paragraph.
Coalesce-Text "()"

# ---------------------------------------------------------------------
# 3) Best-effort: the diff also marks the built-in "Default Paragraph
# Font" character style as semi-hidden (<w:semiHidden/>), matching
# Word's stock template definition for that style. Try the natural COM
# lever for it; harmless no-op if this host doesn't support the
# property.
# ---------------------------------------------------------------------
try {
    $dpf = $d.Styles("Default Paragraph Font")
    $dpf.Hidden = $false
} catch {
}

Write-Output "edit complete"
